$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.755.09'
$ws.Range("E2").Value = '  +0.42%  '

$ws.Range("D3").Value = '1.849.60'
$ws.Range("E3").Value = '  +0.20%  '

$ws.Range("E4").Value = '  -0.09%  '

$ws.Range("D5").Value = '314.21'
$ws.Range("E5").Value = '  -0.15%  '

$ws.Range("E6").Value = '  +0.00%  '

$ws.Range("D7").Value = '0.4337'
$ws.Range("E7").Value = '  +1.60%  '

$ws.Range("D8").Value = '0.3656'
$ws.Range("E8").Value = '  -0.09%  '

$ws.Range("D9").Value = '44.99'
$ws.Range("E9").Value = '  -2.25%  '

$ws.Range("D10").Value = '0.07346'
$ws.Range("E10").Value = '  +1.29%  '

$ws.Range("D11").Value = '0.8781'
$ws.Range("E11").Value = '  -2.46%  '

$ws.Range("D12").Value = '20.77'
$ws.Range("E12").Value = '  +0.38%  '

$ws.Range("D13").Value = '1.821.32'
$ws.Range("E13").Value = '  -1.88%  '

$ws.Range("D14").Value = '5.348'
$ws.Range("E14").Value = '  -0.22%  '

$ws.Range("D15").Value = '6.525'
$ws.Range("E15").Value = '  -0.74%  '

$ws.Range("D16").Value = '0.06958'
$ws.Range("E16").Value = '  +1.82%  '

$ws.Range("D17").Value = '1.003'
$ws.Range("E17").Value = '  -0.02%  '

$ws.Range("D18").Value = '80.12'
$ws.Range("E18").Value = '  +2.90%  '

$ws.Range("D19").Value = '0.000009030'
$ws.Range("E19").Value = '  +2.45%  '

$ws.Range("D20").Value = '1.001'
$ws.Range("E20").Value = '  -0.01%  '

$ws.Range("D21").Value = '15.38'
$ws.Range("E21").Value = '  -0.44%  '

$ws.Range("D22").Value = '27.673.05'
$ws.Range("E22").Value = '  +0.18%  '

$ws.Range("D23").Value = '4.976'
$ws.Range("E23").Value = '  +0.34%  '

$ws.Range("E24").Value = '  -2.29%  '

$ws.Range("D25").Value = '2.036.17'
$ws.Range("E25").Value = '  -2.38%  '

$ws.Range("D26").Value = '1.986'
$ws.Range("E26").Value = '  -2.99%  '

$ws.Range("D27").Value = '155.67'
$ws.Range("E27").Value = '  +1.09%  '

$ws.Range("D28").Value = '18.65'
$ws.Range("E28").Value = '  +2.32%  '

$ws.Range("D29").Value = '120.85'
$ws.Range("E29").Value = '  +9.07%  '

$ws.Range("D30").Value = '5.258'
$ws.Range("E30").Value = '  -0.57%  '

$ws.Range("D31").Value = '1.864'
$ws.Range("E31").Value = '  +2.43%  '

$ws.Range("D32").Value = '0.08924'
$ws.Range("E32").Value = '  +0.45%  '

$ws.Range("D33").Value = '0.7560'
$ws.Range("E33").Value = '  -1.82%  '

$ws.Range("E34").Value = '  +0.00%  '

$ws.Range("D35").Value = '2.974'
$ws.Range("E35").Value = '  +0.00%  '

$ws.Range("D36").Value = '1.123'
$ws.Range("E36").Value = '  +3.56%  '

$ws.Range("D37").Value = '1.108'
$ws.Range("E37").Value = '  +0.81%  '

$ws.Range("D38").Value = '0.05415'
$ws.Range("E38").Value = '  +0.48%  '

$ws.Range("E39").Value = '  +0.47%  '

$ws.Range("D40").Value = '2.825'
$ws.Range("E40").Value = '  -4.40%  '

$ws.Range("D41").Value = '0.5095'
$ws.Range("E41").Value = '  +0.61%  '

$ws.Range("D42").Value = '0.1658'
$ws.Range("E42").Value = '  +0.91%  '

$ws.Range("D43").Value = '6.656'
$ws.Range("E43").Value = '  -2.21%  '

$ws.Range("D44").Value = '8.333'
$ws.Range("E44").Value = '  +1.23%  '

$ws.Range("E45").Value = '  +0.82%  '

$ws.Range("D46").Value = '0.06541'
$ws.Range("E46").Value = '  -1.46%  '

$ws.Range("D47").Value = '0.4667'
$ws.Range("E47").Value = '  -1.22%  '

$ws.Range("D48").Value = '104.36'
$ws.Range("E48").Value = '  -0.87%  '

$ws.Range("E49").Value = '  -0.01%  '

$ws.Range("E50").Value = '  -0.92%  '

$ws.Range("D51").Value = '64.31'
$ws.Range("E51").Value = '  -0.01%  '
